$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.880.47"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.811.54"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'309.69"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4626"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("D9").Value = "'0.07364"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'0.8760"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "'20.47"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "1.837.53"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'5.364"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'6.512"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "'0.07050"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'0.000008705"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'14.74"
$ws.Range("D21").Value = "26.871.78"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'5.312"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "1.970.94"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").Value = "'1.899"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").Value = "'151.53"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "'2.154"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").Value = "'5.321"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'115.98"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "'0.08899"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'0.7540"
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("D33").Value = "'1.156"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.464"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.921"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'1.101"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'0.01966"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "'2.455"
$ws.Range("E39").Value = "  +3.91%  "
$ws.Range("D40").Value = "'0.05260"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'2.917"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'0.5318"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'7.158"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "'0.1661"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "'8.463"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "'0.4959"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'10.28"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.673"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'103.34"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'0.06291"
$ws.Range("E51").Value = "  -1.43%  "
